$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summen")

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 3).Value = 2524081.04254835
    $ws.Cells.Item($r, 6).Value = 29.35664126266641
}
